$d = $word.ActiveDocument

# Locate the bold "*Gestionar Asignatura: ..." paragraph that was appended
# at the very end of the document (after the summary table) and remove it
# entirely, including its paragraph mark, leaving the preceding blank
# (single-space) paragraph intact right before the section break.
$count = $d.Paragraphs.Count
$target = $d.Paragraphs($count)

while ($target.Range.Text.IndexOf("Gestionar Asignatura") -lt 0 -and $count -gt 1) {
    $count = $count - 1
    $target = $d.Paragraphs($count)
}

$target.Range.Delete()
